$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.060.86'
$ws.Range('E2').Value = '  -1.04%  '
$ws.Range('D3').Value = '3.265.97'
$ws.Range('E3').Value = '  -0.45%  '
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').Value = "'583.34"
$ws.Range('E5').Value = '  +0.07%  '
$ws.Range('D6').Value = "'184.86"
$ws.Range('E6').Value = '  -0.48%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = "'0.600"
$ws.Range('E8').Value = '  +0.07%  '
$ws.Range('E9').Value = '  -3.40%  '
$ws.Range('E10').Value = '  -0.92%  '
$ws.Range('E11').Value = '  -3.25%  '
$ws.Range('D12').Value = '3.832.23'
$ws.Range('E12').Value = '  -0.67%  '
$ws.Range('E13').Value = '  +1.42%  '
$ws.Range('D14').Value = "'27.38"
$ws.Range('E14').Value = '  -3.76%  '
$ws.Range('D15').Value = '68.052.90'
$ws.Range('E15').Value = '  -1.10%  '
$ws.Range('D16').Value = "'0.0000167"
$ws.Range('E16').Value = '  -2.44%  '
$ws.Range('D17').Value = '3.282.80'
$ws.Range('E17').Value = '  +0.40%  '
$ws.Range('E18').Value = '  -2.55%  '
$ws.Range('D19').Value = "'13.28"
$ws.Range('E19').Value = '  -2.67%  '
$ws.Range('D20').Value = "'416.24"
$ws.Range('E20').Value = '  +5.23%  '
$ws.Range('E21').Value = '  -2.41%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').Value = "'71.14"
$ws.Range('E23').Value = '  -0.81%  '
$ws.Range('E24').Value = '  -2.41%  '
$ws.Range('E25').Value = '  -3.22%  '
$ws.Range('D26').Value = "'0.187"
$ws.Range('E26').Value = '  -1.14%  '
$ws.Range('E27').Value = '  -3.63%  '
$ws.Range('E28').Value = '  +0.49%  '
$ws.Range('D29').Value = "'1.94"
$ws.Range('E29').Value = '  -2.19%  '
$ws.Range('D30').Value = "'22.61"
$ws.Range('E30').Value = '  -2.37%  '
$ws.Range('E31').Value = '  -4.62%  '
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('E33').Value = '  -4.69%  '
$ws.Range('D34').Value = "'163.40"
$ws.Range('E34').Value = '  -0.04%  '
$ws.Range('E35').Value = '  -5.21%  '
$ws.Range('E36').Value = '  -4.59%  '
$ws.Range('D37').Value = "'26.84"
$ws.Range('E37').Value = '  -0.47%  '
$ws.Range('E38').Value = '  -4.05%  '
$ws.Range('E39').Value = '  -3.92%  '
$ws.Range('D40').Value = "'6.30"
$ws.Range('E40').Value = '  -5.21%  '
$ws.Range('D41').Value = '2.637.75'
$ws.Range('E41').Value = '  -0.83%  '
$ws.Range('E42').Value = '  -4.90%  '
$ws.Range('E43').Value = '  -2.80%  '
$ws.Range('D44').Value = "'336.94"
$ws.Range('E44').Value = '  -1.86%  '
$ws.Range('D45').Value = "'24.23"
$ws.Range('E45').Value = '  -5.00%  '
$ws.Range('E46').Value = '  -3.70%  '
$ws.Range('E47').Value = '  -2.15%  '
$ws.Range('E48').Value = '  -2.08%  '
$ws.Range('E49').Value = '  -2.09%  '
$ws.Range('E50').Value = '  -0.09%  '
$ws.Range('D51').Value = "'30.52"
$ws.Range('E51').Value = '  -5.11%  '

foreach ($addr in @('D5','D6','D8','D14','D16','D19','D20','D23','D26','D29','D30','D34','D37','D40','D44','D45','D51')) {
    $ws.Range($addr).Style = "Normal"
}
